# Updates the cryptos price/volume table (columns D and E) on the active
# worksheet to the latest scraped values, matching the upstream GitHub
# Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.169.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.51%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.780.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.67%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.01%  '

# Row 7
$ws.Range("E7").Value = '  +1.90%  '

# Row 8
$ws.Range("E8").Value = '  +0.27%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.790.38'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.98%  '

# Row 10
$ws.Range("E10").Value = '  +0.89%  '

# Row 11
$ws.Range("E11").Value = '  +1.55%  '

# Row 12
$ws.Range("E12").Value = '  +3.22%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.160'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.15%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.275.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.62%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.55%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.101.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.59%  '

# Row 17
$ws.Range("E17").Value = '  +6.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.787.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.40%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.02%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '367.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.34%  '

# Row 22
$ws.Range("E22").Value = '  +0.87%  '

# Row 23
$ws.Range("E23").Value = '  +8.05%  '

# Row 24
$ws.Range("E24").Value = '  +0.67%  '

# Row 25
$ws.Range("E25").Value = '  +3.36%  '

# Row 26
$ws.Range("E26").Value = '  +6.50%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.14%  '

# Row 28
$ws.Range("E28").Value = '  +12.59%  '

# Row 29
$ws.Range("E29").Value = '  +0.50%  '

# Row 30
$ws.Range("E30").Value = '  +1.26%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.29'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.02%  '

# Row 32
$ws.Range("E32").Value = '  +5.19%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '172.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.24%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.65%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.89'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.68%  '

# Row 36
$ws.Range("E36").Value = '  +0.11%  '

# Row 37
$ws.Range("E37").Value = '  +6.28%  '

# Row 38
$ws.Range("E38").Value = '  +2.07%  '

# Row 39
$ws.Range("E39").Value = '  +2.65%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '343.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.08%  '

# Row 41
$ws.Range("E41").Value = '  +0.44%  '

# Row 42
$ws.Range("E42").Value = '  +11.21%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.82%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.34%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.27%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0612'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.11%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.654'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.05%  '

# Row 48
$ws.Range("E48").Value = '  +1.23%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.27%  '

# Row 50
$ws.Range("E50").Value = '  +2.34%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.179.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.66%  '
